$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Monday 23.4.18, 1200-1500, 3 hours, Theory + interface...
$ws.Cells.Item(7, 1).Value = "Monday 23.4.18"
$ws.Cells.Item(7, 2).Value = "1200-1500"
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = "Theory + interface (ComparePolygons) + subclass (SortedPolygons)"

# Row 8: Monday 23.4.18, 1630-1730, 1 hour, Getting comparePolygons to work
$ws.Cells.Item(8, 1).Value = "Monday 23.4.18"
$ws.Cells.Item(8, 2).Value = "1630-1730"
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = "Getting comparePolygons to work"

# Update selection to B9
$ws.Range("B9").Select()
